$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: keep the existing thin-box border but drop left & right -> top+bottom only
$ws1.Range("C1").Borders(7).LineStyle = -4142   # xlEdgeLeft  -> none
$ws1.Range("C1").Borders(10).LineStyle = -4142  # xlEdgeRight -> none

# D1: keep the existing thin-box border but drop left -> top+bottom+right
$ws1.Range("D1").Borders(7).LineStyle = -4142   # xlEdgeLeft -> none

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Borders(7).LineStyle = -4142
$ws2.Range("C1").Borders(10).LineStyle = -4142

$ws2.Range("D1").Borders(7).LineStyle = -4142

$ws2.Range("F1").Borders(7).LineStyle = -4142
$ws2.Range("F1").Borders(10).LineStyle = -4142

$ws2.Range("G1").Borders(7).LineStyle = -4142

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 previously held an empty inline string cell; drop it entirely
$ws2.Range("G5").ClearContents()
